$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column BH: header "Agosto.2021" following "Mayo.2021" in BG1
$ws.Range("BH1").Value = "Agosto.2021"

# Copy formatting from BG1 (last header cell) to BH1 so it matches the other headers
$ws.Range("BG1").Copy()
$ws.Range("BH1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill column BH (rows 2-19) with the same value as column BG for each row
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 60).Formula = $ws.Cells.Item($r, 59).Formula
}
